$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44253
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("S2").Value = 806

# Row 3
$ws.Range("D3").Value = 45072
$ws.Range("L3").Value = "Segunda"
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 16000
$ws.Range("R3").Value = "Provincia de Chacabuco"
$ws.Range("S3").Value = 889

# Row 4
$ws.Range("D4").Value = 45072
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 17000
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 944

# Row 5
$ws.Range("D5").Value = 44250
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 806

# Row 6
$ws.Range("D6").Value = 44252
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 14000
$ws.Range("P6").Value = 13500
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("S6").Value = 750

# Row 7
$ws.Range("D7").Value = 44257
$ws.Range("M7").Value = 100

